# Sync attendance_reports: fix "Recorded By" (column G) ordering.
# Rule: values are comma-separated lists (e.g. "System, someone@example.com").
# Whenever the FIRST item in the list is "System" (case-insensitive), it is
# moved to swap places with the LAST item in the list, so "System" ends up
# trailing instead of leading (e.g. "System, a@b.com" -> "a@b.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G is the "Recorded By" column.
$col = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        if ($parts.Count -ge 2 -and $parts[0].ToLower() -eq "system") {
            $first = $parts[0]
            $last = $parts[$parts.Count - 1]
            $parts[0] = $last
            $parts[$parts.Count - 1] = $first
            $newVal = [string]::Join(", ", $parts)
            $cell.Value = $newVal
        }
    }
}
